$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing data rows (e.g. row 61) carry a center/center-aligned style.
# Copy that formatting onto the new row 62 first so the new cells match
# the style used by the rest of the table.
$ws.Range("A61:C61").Copy()
$ws.Range("A62:C62").PasteSpecial(-4122)

# Column A stores the date as literal text (e.g. "2026/01/10"), not as a
# real date value. Temporarily mark A62 as Text before typing the value so
# Excel doesn't auto-convert the "2026/01/11" string into a date serial.
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "2026/01/11"

# Re-apply the row formatting so A62 ends up with the same style as the
# rest of the column (General number format, centered alignment) instead
# of staying tagged as Text.
$ws.Range("A61:C61").Copy()
$ws.Range("A62:C62").PasteSpecial(-4122)

$ws.Range("B62").Value = "逃离鸭科夫"
$ws.Range("C62").Value = 1143
